$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.942.55"
$ws.Range("E2").Value = "  -1.49%  "
$ws.Range("D3").Value = "3.097.95"
$ws.Range("E3").Value = "  -2.48%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.28%  "
$ws.Range("D5").Value = "'593.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.22%  "
$ws.Range("D6").Value = "'156.60"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.23%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "3.098.14"
$ws.Range("E9").Value = "  -2.48%  "
$ws.Range("E10").Value = "  -3.57%  "
$ws.Range("D11").Value = "'5.92"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.50%  "
$ws.Range("E12").Value = "  -3.68%  "
$ws.Range("E13").Value = "  -5.04%  "
$ws.Range("D14").Value = "'37.00"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.40%  "
$ws.Range("E15").Value = "  -1.26%  "
$ws.Range("D16").Value = "3.607.53"
$ws.Range("E16").Value = "  -2.62%  "
$ws.Range("E17").Value = "  -1.93%  "
$ws.Range("D18").Value = "63.798.70"
$ws.Range("E18").Value = "  -1.35%  "
$ws.Range("D19").Value = "3.098.02"
$ws.Range("E19").Value = "  -2.71%  "
$ws.Range("D20").Value = "'479.85"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.09%  "
$ws.Range("D21").Value = "'14.49"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.31%  "
$ws.Range("D22").Value = "'0.713"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.69%  "
$ws.Range("D23").Value = "'7.57"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.76%  "
$ws.Range("D24").Value = "'2.47"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.58%  "
$ws.Range("D25").Value = "'81.31"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.48%  "
$ws.Range("D26").Value = "'12.92"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.53%  "
$ws.Range("D27").Value = "'10.67"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +6.82%  "
$ws.Range("E28").Value = "  -0.16%  "
$ws.Range("D29").Value = "'7.53"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.87%  "
$ws.Range("E30").Value = "  -2.58%  "
$ws.Range("D31").Value = "'0.999"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.42%  "
$ws.Range("D32").Value = "'2.19"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.01%  "
$ws.Range("E33").Value = "  -4.72%  "
$ws.Range("D34").Value = "'27.27"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.60%  "
$ws.Range("D35").Value = "0.0₃0842"
$ws.Range("E35").Value = "  -6.14%  "
$ws.Range("E36").Value = "  -2.37%  "
$ws.Range("D37").Value = "'6.03"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.20%  "
$ws.Range("D38").Value = "'2.26"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.04%  "
$ws.Range("D39").Value = "'3.29"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -7.86%  "
$ws.Range("D40").Value = "'50.95"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.82%  "
$ws.Range("D41").Value = "'9.23"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.83%  "
$ws.Range("D42").Value = "'440.90"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.52%  "
$ws.Range("D43").Value = "'0.291"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.51%  "
$ws.Range("D44").Value = "'0.0365"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.96%  "
$ws.Range("E45").Value = "  +0.82%  "
$ws.Range("D46").Value = "'40.14"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.80%  "
$ws.Range("D47").Value = "2.828.71"
$ws.Range("E47").Value = "  -3.61%  "
$ws.Range("D48").Value = "'131.27"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.06%  "
$ws.Range("D49").Value = "'25.94"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.09%  "
$ws.Range("E50").Value = "  +0.02%  "
$ws.Range("E51").Value = "  -3.51%  "
